# Corrected stellar parameters retrieval:
# - Prefix the Gaia source id in column A with "Gaia DR3" (or "Gaia DR2"
#   for the rows whose DR3 id, column C, is blank/unavailable).
# - Widen column A so the longer labels remain fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $idCell = $ws.Cells.Item($r, 1)
    $idText = $idCell.Text

    if ($idText -eq "") {
        continue
    }

    $dr3Cell = $ws.Cells.Item($r, 3)

    if ($dr3Cell.Text -eq "") {
        $prefix = "Gaia DR2"
    } else {
        $prefix = "Gaia DR3"
    }

    $idCell.Value = "$prefix $idText"
}

# Raw OOXML column width of 30 <=> COM ColumnWidth of 29.17 (offset of
# 0.8333... between the two representations for this font/theme).
$ws.Columns.Item(1).ColumnWidth = 29.17
